$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") rows 2-18: update the serial date value from 45192 to 45202
for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 3).Value = 45202
}
